$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the execution time for "js in browser" (row 2) from 1h 19m to 1h 55m
$ws.Range("C2").Value = "1h 55m"

# Reflect the active cell selection recorded in the saved file (C3)
$ws.Range("C3").Select()
